$wb = $excel.ActiveWorkbook

# Update the value of C3 on the "Data" sheet: "chrome" -> "chromegrid"
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("C3").Value = "chromegrid"

# Set selection on Data sheet to C3 (matches final saved state before switching tabs)
$wsData.Range("C3").Select()

# Switch active sheet to "Test" and select C3 there
$wsTest = $wb.Worksheets.Item("Test")
$wsTest.Activate()
$wsTest.Range("C3").Select()
